# Apply edits described by the commit "fg wb xls update"
$wb = $excel.ActiveWorkbook

# --- Sheet: cost calculation ---
$wsCost = $wb.Worksheets.Item("cost calculation")

# One more "forest goblins" henchmen unit bought (quantity 17 -> 18)
$wsCost.Range("E14").Value = 18

# Manually updated "Used cost" figure next to the "+1" Warband Equip Rating row
$wsCost.Range("K13").Value = 209

# Simplify/fix the "Warband Rating" formula to reference E18 (total units) directly
# instead of routing through J3 and the (always zero) Hired swords!B5 term
$wsCost.Range("J2").Formula = "=E18+H38+I38+J38"

$wsCost.Activate()
$wsCost.Range("K14").Select()

# --- Sheet: Underdog Bonus ---
$wsUnderdog = $wb.Worksheets.Item("Underdog Bonus")

# Update the underdog-bonus lookup row formula to use the new total units (E18)
# and add the hired-swords upkeep base (Hired swords!B19) instead of the old
# (J3 - E17) difference
$wsUnderdog.Range("C11:U11").Formula = "=C4*'cost calculation'!`$E`$18+'Hired swords'!`$B`$19"

$wsUnderdog.Activate()
$wsUnderdog.Range("L15").Select()

# --- Sheet: Hired swords ---
$wsHired = $wb.Worksheets.Item("Hired swords")
$wsHired.Activate()
$wsHired.Range("B19").Select()

# Restore the originally active/selected sheet ("cost calculation" stays the
# tab that is shown/selected when the workbook is opened)
$wsCost.Activate()
